$wb = $excel.ActiveWorkbook

# Sheet names affected: "展览" and "全部类型" (both contain duplicated event data)
$sheetNames = @("展览", "全部类型")

# New values for column F ("想去人数") keyed by row number
$updates = @{
    2 = 1333
    3 = 1863
    4 = 156
    6 = 6307
    7 = 173
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
